$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") contains a date serial number (45171 -> 2023-09-02)
# that was bumped by one day (45172 -> 2023-09-03) for every data row.
$ws.Range("C2:C454").Value = 45172
